$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected cells are formatted as Text so that numeric-looking
# strings (e.g. "0.06270", "16.10") keep their exact original formatting
# (significant digits, separators, etc.) instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.950.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5236"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2613"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06270"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07733"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.658.58"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.447"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5427"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8054"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.992.24"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.963"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1236"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.248"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05922"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.273"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.485"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.228"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.519"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.416"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9396"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5692"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.844"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8459"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.44"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "999.62"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.786.87"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.43"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4292"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.477"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05150"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.817"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.36%  "
